$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Software": remove the now-obsolete "Version/Spec" column (I) and
# update a handful of cell values (Minio chart/app bump, an extra note line,
# and two TODO cleanups).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Software")

# Drop column I ("Version/Spec") - everything to its right shifts left.
$ws.Columns("I").Delete() | Out-Null

# Minio row: bump chart/app versions and correct the memory allocation.
$ws.Range("D9").Value = "minio-2.3.4"
$ws.Range("E9").Value = "minio/minio:RELEASE.2018-12-06T01-27-43Z"
$ws.Range("G9").Value = 0.25

# Prometheus row: append a third line to the existing note.
$note = "1. Change all TCP NodePorts to Nginx Ports`n2. Mongo Exporter (Evelyn)`n3. recycling mechanism (Evelyn)"
$ws.Range("N10").Value = $note

# RMMPortal / RMMWorker / OTA-Worker rows: "TODO (Alex)" -> "TODO"
$ws.Range("E14").Value = "TODO"
$ws.Range("E15").Value = "TODO"
$ws.Range("E16").Value = "TODO"

# This sheet is no longer the active tab; remember a different selection.
$ws.Range("E22").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "UserManual" becomes the active tab with a new selection.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("UserManual")
$ws3.Activate() | Out-Null
$ws3.Range("D6").Select() | Out-Null
